# "added max capacity for all options"
#
# Adds a new "Capacity limit grid" column (with a "-" unit and a 100
# value) to the "General Data" sheet, and replays the view/selection
# state (active sheet + per-sheet last-selected cell) left behind by the
# editing session.

$wb = $excel.ActiveWorkbook

# --- Sets: just visited / a cell selected there, no data changed -----
$wsSets = $wb.Worksheets.Item("Sets")
$wsSets.Activate()
$wsSets.Range("E3").Select() | Out-Null
$wsSets.Columns.Item(4).AutoFit() | Out-Null

# --- irradiation: visited / a cell selected there, no data changed ---
$wsIrr = $wb.Worksheets.Item("irradiation")
$wsIrr.Activate()
$wsIrr.Range("C1").Select() | Out-Null
$wsIrr.Columns.Item(2).AutoFit() | Out-Null
$wsIrr.Columns.Item(3).ColumnWidth = 11.5

# --- General Data: new "Capacity limit grid" column -------------------
$wsGen = $wb.Worksheets.Item("General Data")
$wsGen.Activate()

$wsGen.Range("E1").Value = "Capacity limit grid"
$wsGen.Range("E2").Value = "-"
$wsGen.Range("E3").Value = 100

$wsGen.Columns.Item(5).AutoFit() | Out-Null

$wsGen.Range("E1").Select() | Out-Null
